$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price column D holds numeric-looking values that the source sheet
# stores as plain text (inline strings), not real numbers. A leading
# apostrophe forces Excel to keep each assignment as text instead of
# auto-converting it to a number, matching the workbook's original typing.

$ws.Range("D2").Value = '''244.88'
$ws.Range("D3").Value = '''21.95'
$ws.Range("D4").Value = '''5.394'
$ws.Range("D5").Value = '''0.05983'
$ws.Range("D7").Value = '''6.383'
$ws.Range("D8").Value = '''0.8107'
$ws.Range("D9").Value = '''0.9618'
$ws.Range("D10").Value = '''0.1428'
$ws.Range("D11").Value = '''0.07396'
$ws.Range("D12").Value = '''0.03398'
$ws.Range("D13").Value = '''0.03062'
$ws.Range("D14").Value = '''0.09419'
$ws.Range("D15").Value = '''4.005'
$ws.Range("D16").Value = '''0.001598'
$ws.Range("D17").Value = '''0.04807'
$ws.Range("D18").Value = '''0.0005872'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("D19").Value = '''0.006129'
$ws.Range("D20").Value = '''0.005070'
$ws.Range("D21").Value = '''0.0009822'
$ws.Range("D22").Value = '''0.0001000'
$ws.Range("D23").Value = '''3.726'
$ws.Range("D24").Value = '''2.186'
$ws.Range("D40").Value = '''0.03975'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '''0.006582'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1073'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '''0.002901'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = '''0.005306'
$ws.Range("D45").Value = '''0.00005258'
$ws.Range("D48").Value = '''0.02647'
$ws.Range("E48").Value = '47BOLOBOLO'
